$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: date value change
$ws.Range("A1").Value = 45352

# Row 2 picks up an explicit (custom) row height matching the default 15.05
$ws.Rows(2).RowHeight = 15.05

# Row 4
$ws.Range("A4").Value = 163575
$ws.Range("B4").Value = 100.77244688980591
$ws.Range("C4").Value = 16483853

# Row 5
$ws.Range("A5").Value = 2645
$ws.Range("B5").Value = 236.51228733459357
$ws.Range("C5").Value = 625575

# Row 6
$ws.Range("A6").Value = 46
$ws.Range("B6").Value = 2247.3478260869565
$ws.Range("C6").Value = 103378

# Row 7
$ws.Range("A7").Value = 166266
$ws.Range("B7").Value = 103.52571181119411
$ws.Range("C7").Value = 17212806

# Row 9
$ws.Range("A9").Value = 33018
$ws.Range("B9").Value = 165
$ws.Range("C9").Value = 5447970

# Row 10
$ws.Range("A10").Value = 65057
$ws.Range("B10").Value = 165
$ws.Range("C10").Value = 10734405

# Row 11
$ws.Range("A11").Value = 98075
$ws.Range("B11").Value = 165
$ws.Range("C11").Value = 16182375

# Row 12: the whole row content is cleared (was a shared-string " " placeholder cell in A12)
# Row numbers of following rows stay the same, so use ClearContents (no shifting).
$ws.Range("A12:C12").ClearContents()

# Row 13
$ws.Range("A13").Value = 294635
$ws.Range("B13").Value = 229
$ws.Range("C13").Value = 67471415

# Row 14
$ws.Range("A14").Value = 848
$ws.Range("B14").Value = 218.00943396226415
$ws.Range("C14").Value = 184872

# Row 15: value moves from B15 to A15
$ws.Range("B15").ClearContents()
$ws.Range("A15").Value = 0

# Row 16
$ws.Range("A16").Value = 295483
$ws.Range("B16").Value = 228.96845842231195
$ws.Range("C16").Value = 67656287

# Row 19: cell C19 content removed entirely (row number of row 21/22 unchanged)
$ws.Range("C19").ClearContents()

# Row 21
$ws.Range("A21").Value = 559824
$ws.Range("B21").Value = 180.50578038812199
$ws.Range("C21").Value = 101051468

# Row 22
$ws.Range("A22").Value = 393558
$ws.Range("B22").Value = 213.02746228001971
$ws.Range("C22").Value = 83838662

# Update selection to A2
$ws.Range("A2").Select()
